$d = $word.ActiveDocument

# 1) Insert a new paragraph right after the "admin" paragraph, carrying the
#    same paragraph formatting (no strike), and give it the new sentence.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "service has /hour and + should i change pricing to string"

# 2) Replace the text of the (now third) paragraph that used to read
#    "how to do update product for frontend or call page" while leaving its
#    paragraph/run formatting (including the strike-through) untouched.
$d.Content.Find.Execute("how to do update product for frontend or call page", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ask for seeding the relationship product", 2)
